$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1397.528150389182
$ws.Range("C3").Value = 1551.173121818517
$ws.Range("C4").Value = 1530.527801428039
$ws.Range("C5").Value = 1416.937498208815
$ws.Range("C6").Value = 1411.876087204833
$ws.Range("C7").Value = 1380.959382710767
$ws.Range("C8").Value = 1294.047308086022
$ws.Range("C9").Value = 1341.021254713692
$ws.Range("C10").Value = 1362.879484263613
$ws.Range("C11").Value = 1314.052484074679
$ws.Range("C12").Value = 1305.338538540757
$ws.Range("C13").Value = 1321.586462558755
$ws.Range("C14").Value = 1311.477870204247
$ws.Range("C15").Value = 1306.191350510799
$ws.Range("C16").Value = 1278.14633747203
$ws.Range("C17").Value = 1261.414924626947
$ws.Range("C18").Value = 1230.934146916329
$ws.Range("C19").Value = 1236.652382254011
$ws.Range("C20").Value = 1233.635855764303
$ws.Range("C21").Value = 1223.959023374688
$ws.Range("C22").Value = 1209.43629271808
$ws.Range("C23").Value = 1193.74676171502
$ws.Range("C24").Value = 1178.679177740949
$ws.Range("C25").Value = 1169.296754069848
$ws.Range("C26").Value = 1159.671040694144
$ws.Range("C27").Value = 1150.428482618038
$ws.Range("C28").Value = 1137.020365836294
$ws.Range("C29").Value = 1129.110785367038
$ws.Range("C30").Value = 1123.887709604254
$ws.Range("C31").Value = 1114.153154214592
$ws.Range("C32").Value = 1108.282335600008
$ws.Range("C33").Value = 1102.75255490766
$ws.Range("C34").Value = 1097.806806432677
$ws.Range("C35").Value = 1094.063075973706
$ws.Range("C36").Value = 1090.552839180692
$ws.Range("C37").Value = 1086.466195104219
$ws.Range("C38").Value = 1084.363847458136
$ws.Range("C39").Value = 1080.906383514768
$ws.Range("C40").Value = 1076.105202275746
$ws.Range("C41").Value = 1072.056951943679
$ws.Range("C42").Value = 1071.620925971684
$ws.Range("C43").Value = 1068.428865162341
$ws.Range("C44").Value = 1066.666852495235
$ws.Range("C45").Value = 1064.963983419043
$ws.Range("C46").Value = 1064.400579256497
$ws.Range("C47").Value = 1063.02604946884
$ws.Range("C48").Value = 1060.639142081765
$ws.Range("C49").Value = 1060.038348490666
$ws.Range("C50").Value = 1057.840267100985
$ws.Range("C51").Value = 1054.615982252871
$ws.Range("C52").Value = 1051.900917514836
$ws.Range("C53").Value = 1045.836962149216
$ws.Range("C54").Value = 1040.699890698216
$ws.Range("C55").Value = 1032.318108457608
$ws.Range("C56").Value = 1023.953178870718
$ws.Range("C57").Value = 1020.116562746809
$ws.Range("C58").Value = 1018.778608112614
$ws.Range("C59").Value = 1016.679883241167
$ws.Range("C60").Value = 1009.418890728808
$ws.Range("C61").Value = 993.0143088973109
$ws.Range("C62").Value = 992.7675703351734
$ws.Range("C63").Value = 988.5348692387774
$ws.Range("C64").Value = 987.222952048561
$ws.Range("C65").Value = 987.9894054248828
